$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two "ready rack" ammo-loading helper rows (e009a / e009b) are being
# replaced by a single new "e010 Time Check" row. Remove the old rows 11-12
# and insert a fresh row 11 in their place (this keeps row 10 and all rows
# below row 12 shifting up naturally, matching the target layout).
$ws.Rows("11:12").Delete()
$ws.Rows("11").Insert()

# Populate the new row 11 with the e010 event id + its body text.
$ws.Range("A11").Value = "e010"

$e010Body = "<Bold>e010 Time Check</Bold> `n" + `
  "<InlineUIContainer><Button Content='r4.3' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  `n" + `
  "<InlineUIContainer><Button Content='Time' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>`n" + `
  "<LineBreak/><LineBreak/>`n" + `
  "Determine sunrise and sunset for current month using the <InlineUIContainer><Button Content='Time' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Table. The same die roll is used to determine the ammo expended:`n" + `
  "<InlineUIContainer><Image Name='DieRoll' Height='21' Width='21' > </Image></InlineUIContainer>."

$ws.Range("B11").Value = $e010Body

# Match the row height Excel computed for the new wrapped text.
$ws.Rows("11").RowHeight = 99.85

# Update the selection to match where the author left off editing.
$ws.Range("B13:B14").Select()
